$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.343.77"
$ws.Range("E2").Value = "  +0.17%  "

$ws.Range("D3").Value = "1.798.59"
$ws.Range("E3").Value = "  +0.07%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").Value = "'227.00"
$ws.Range("E5").Value = "  +0.01%  "

$ws.Range("E6").Value = "  +3.07%  "

$ws.Range("E7").Value = "  -0.03%  "

$ws.Range("D8").Value = "'36.19"
$ws.Range("E8").Value = "  +9.61%  "

$ws.Range("D9").Value = "'0.299"
$ws.Range("E9").Value = "  +1.43%  "

$ws.Range("E10").Value = "  +0.22%  "

$ws.Range("D11").Value = "'0.0961"
$ws.Range("E11").Value = "  +1.61%  "

$ws.Range("D12").Value = "2.056.53"
$ws.Range("E12").Value = "  +0.01%  "

$ws.Range("D13").Value = "'11.53"
$ws.Range("E13").Value = "  +4.28%  "

$ws.Range("D14").Value = "1.794.52"
$ws.Range("E14").Value = "  -15.65%  "

$ws.Range("D15").Value = "'0.640"
$ws.Range("E15").Value = "  +1.24%  "

$ws.Range("D16").Value = "'4.48"
$ws.Range("E16").Value = "  +4.99%  "

$ws.Range("D17").Value = "34.295.89"
$ws.Range("E17").Value = "  -0.03%  "

$ws.Range("D18").Value = "'68.78"
$ws.Range("E18").Value = "  +0.52%  "

$ws.Range("D19").Value = "'243.81"
$ws.Range("E19").Value = "  +0.05%  "

$ws.Range("D20").Value = "0.0₃0789"
$ws.Range("E20").Value = "  -0.47%  "

$ws.Range("D21").Value = "'11.58"
$ws.Range("E21").Value = "  +2.86%  "

$ws.Range("E22").Value = "  +0.07%  "

$ws.Range("D23").Value = "'4.16"
$ws.Range("E23").Value = "  +0.06%  "

$ws.Range("E24").Value = "  +3.87%  "

$ws.Range("D25").Value = "'172.02"
$ws.Range("E25").Value = "  +3.28%  "

$ws.Range("D26").Value = "'7.94"
$ws.Range("E26").Value = "  +8.75%  "

$ws.Range("D27").Value = "'16.77"
$ws.Range("E27").Value = "  +1.81%  "

$ws.Range("E28").Value = "  +1.30%  "

$ws.Range("E29").Value = "  -0.13%  "

$ws.Range("D30").Value = "'3.99"
$ws.Range("E30").Value = "  +0.95%  "

$ws.Range("E31").Value = "  +0.10%  "

$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value = "'3.82"
$ws.Range("E32").Value = "  +1.12%  "

$ws.Range("B33").Value = "PancakeSwap"
$ws.Range("C33").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D33").Value = "'1.24"
$ws.Range("E33").Value = "  +0.27%  "

$ws.Range("D34").Value = "'1.81"
$ws.Range("E34").Value = "  -0.61%  "

$ws.Range("D35").Value = "1.394.58"
$ws.Range("E35").Value = "  -0.51%  "

$ws.Range("D36").Value = "'0.669"
$ws.Range("E36").Value = "  +0.21%  "

$ws.Range("D37").Value = "'2.45"
$ws.Range("E37").Value = "  -6.06%  "

$ws.Range("E38").Value = "  -0.34%  "

$ws.Range("E39").Value = "  +0.03%  "

$ws.Range("D40").Value = "'0.957"
$ws.Range("E40").Value = "  +2.29%  "

$ws.Range("D41").Value = "'82.24"

$ws.Range("E42").Value = "  -0.72%  "

$ws.Range("D43").Value = "'2.42"
$ws.Range("E43").Value = "  +0.18%  "

$ws.Range("D44").Value = "'1.19"
$ws.Range("E44").Value = "  +6.54%  "

$ws.Range("D45").Value = "'13.33"
$ws.Range("E45").Value = "  -4.25%  "

$ws.Range("D46").Value = "'6.00"
$ws.Range("E46").Value = "  -0.33%  "

$ws.Range("E47").Value = "  -4.14%  "

$ws.Range("D48").Value = "1.957.51"
$ws.Range("E48").Value = "  +0.07%  "

$ws.Range("E49").Value = "  +0.01%  "

$ws.Range("D50").Value = "'103.89"
$ws.Range("E50").Value = "  -0.85%  "

$ws.Range("D51").Value = "0.0₆0127"
$ws.Range("E51").Value = "  -0.89%  "
